# Updated cryptos list on Thu May 16 19:58:42 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving the cell's original style
# (prevents Excel from auto-coercing numeric-looking strings like "6.68"
# or "65.253.12" into numbers).
function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

# --- Row 2: Bitcoin ---
Set-TextValue 2 4 "65.253.12"
$ws.Range("E2").Value = "  -1.20%  "

# --- Row 3: Ethereum ---
Set-TextValue 3 4 "2.946.57"
$ws.Range("E3").Value = "  -2.34%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.02%  "

# --- Row 5: BNB ---
Set-TextValue 5 4 "567.99"
$ws.Range("E5").Value = "  -3.12%  "

# --- Row 6: Solana ---
Set-TextValue 6 4 "159.26"
$ws.Range("E6").Value = "  +3.01%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  +0.03%  "

# --- Row 8: XRP ---
$ws.Range("E8").Value = "  -0.04%  "

# --- Row 9: LidoStakedEther ---
Set-TextValue 9 4 "2.942.09"
$ws.Range("E9").Value = "  -2.39%  "

# --- Row 10: Toncoin ---
Set-TextValue 10 4 "6.68"
$ws.Range("E10").Value = "  -4.35%  "

# --- Row 11: Dogecoin ---
$ws.Range("E11").Value = "  -2.39%  "

# --- Row 12: Cardano ---
Set-TextValue 12 4 "0.458"
$ws.Range("E12").Value = "  +1.43%  "

# --- Row 13: ShibaInu ---
Set-TextValue 13 4 "0.0000244"
$ws.Range("E13").Value = "  +0.98%  "

# --- Row 14: Avalanche ---
$ws.Range("E14").Value = "  +0.22%  "

# --- Row 15: TRON ---
$ws.Range("E15").Value = "  -0.74%  "

# --- Row 16: WrappedBTC ---
Set-TextValue 16 4 "65.297.70"
$ws.Range("E16").Value = "  -1.09%  "

# --- Row 17: WrappedliquidstakedEther2.0 ---
Set-TextValue 17 4 "3.438.64"
$ws.Range("E17").Value = "  -2.22%  "

# --- Row 18: Polkadot ---
Set-TextValue 18 4 "6.96"
$ws.Range("E18").Value = "  -0.31%  "

# --- Row 19: WrappedEther ---
Set-TextValue 19 4 "2.951.39"
$ws.Range("E19").Value = "  -2.08%  "

# --- Row 20: Chainlink ---
Set-TextValue 20 4 "14.85"
$ws.Range("E20").Value = "  +7.20%  "

# --- Row 21: BitcoinCash ---
Set-TextValue 21 4 "445.32"
$ws.Range("E21").Value = "  -2.73%  "

# --- Row 22: Polygon ---
Set-TextValue 22 4 "0.687"
$ws.Range("E22").Value = "  +0.12%  "

# --- Row 23: Uniswap ---
Set-TextValue 23 4 "7.23"
$ws.Range("E23").Value = "  -1.84%  "

# --- Row 24: Litecoin ---
Set-TextValue 24 4 "82.20"
$ws.Range("E24").Value = "  +0.44%  "

# --- Row 25: Fetch.AI ---
$ws.Range("E25").Value = "  -2.05%  "

# --- Row 26: InternetComputer(DFINITY) ---
Set-TextValue 26 4 "12.11"
$ws.Range("E26").Value = "  -4.33%  "

# --- Rows 27 & 28: RenderToken / Dai swap places in the ranking ---
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue 27 4 "1.00"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 28 4 "10.00"
$ws.Range("E28").Value = "  -6.79%  "

# --- Row 29: NEARProtocol ---
$ws.Range("E29").Value = "  +1.51%  "

# --- Row 30: ImmutableX ---
$ws.Range("E30").Value = "  -2.15%  "

# --- Row 31: PancakeSwap ---
$ws.Range("E31").Value = "  -1.29%  "

# --- Row 32: PEPE ---
$ws.Range("E32").Value = "  -2.50%  "

# --- Row 33: EthereumClassic ---
Set-TextValue 33 4 "27.19"
$ws.Range("E33").Value = "  +0.52%  "

# --- Row 34: Hedera ---
Set-TextValue 34 4 "0.110"
$ws.Range("E34").Value = "  -1.29%  "

# --- Row 35: FirstDigitalUSD ---
Set-TextValue 35 4 "1.00"
$ws.Range("E35").Value = "  +0.05%  "

# --- Row 36: Mantle ---
Set-TextValue 36 4 "0.975"
$ws.Range("E36").Value = "  -1.56%  "

# --- Row 37: Filecoin ---
Set-TextValue 37 4 "5.75"
$ws.Range("E37").Value = "  -0.69%  "

# --- Row 38: OKB ---
Set-TextValue 38 4 "49.33"
$ws.Range("E38").Value = "  -0.10%  "

# --- Row 39: Arweave ---
$ws.Range("E39").Value = "  -2.70%  "

# --- Row 40: Stacks ---
$ws.Range("E40").Value = "  -7.99%  "

# --- Rows 41 & 42: dogwifhat / Kaspa swap places in the ranking ---
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue 41 4 "0.120"
$ws.Range("E41").Value = "  -1.40%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue 42 4 "2.83"
$ws.Range("E42").Value = "  -4.00%  "

# --- Row 43: TheGraph ---
$ws.Range("E43").Value = "  -1.58%  "

# --- Row 44: Cosmos ---
Set-TextValue 44 4 "8.46"
$ws.Range("E44").Value = "  -0.07%  "

# --- Row 45: Bittensor ---
Set-TextValue 45 4 "384.57"
$ws.Range("E45").Value = "  -0.81%  "

# --- Row 46: VeChain ---
$ws.Range("E46").Value = "  -0.79%  "

# --- Row 47: Maker ---
Set-TextValue 47 4 "2.706.87"
$ws.Range("E47").Value = "  -2.82%  "

# --- Row 48: Monero ---
Set-TextValue 48 4 "132.93"
$ws.Range("E48").Value = "  -1.38%  "

# --- Row 49: USDe ---
$ws.Range("E49").Value = "  +0.05%  "

# --- Row 50: ThetaToken ---
$ws.Range("E50").Value = "  +4.44%  "

# --- Row 51: Stellar ---
$ws.Range("E51").Value = "  -0.10%  "
